# "Ajout de la colonne Archétypes"
# Insert a new "Archétype" column into the "Espèces" sheet (between the
# existing "Concept" column C and the old "Nb Tribus" column, which was D
# and becomes E), then fill in the archetype for each species that has one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Espèces")

# Shift columns D:L one place to the right to make room for the new column.
$ws.Range("D1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("D1").Value = "Archétype"

# Species archétypes (rows follow the existing species order).
# Values are entered in the same order the original author typed them in
# (this matters only for shared-string ordering, not for the final data).
$ws.Range("D3").Value  = "Comploteur"    # Ananasi
$ws.Range("D4").Value  = "Érudit"        # Bastet
$ws.Range("D5").Value  = "Messager"      # Corax
$ws.Range("D6").Value  = "Guerrisseur"   # Gurahl
$ws.Range("D8").Value  = "Mémoire"       # Mokolé
$ws.Range("D9").Value  = "Assassin"      # Nâga
$ws.Range("D10").Value = "Farceur"       # Nuwisha
$ws.Range("D13").Value = "Traqueur"      # Uratha
$ws.Range("D7").Value  = "Manipulateur"  # Kitsune

# Re-apply the AutoFilter over the new, wider range (A1:M13).
$ws.AutoFilterMode = $false
$ws.Range("A1:M13").AutoFilter()

# Update the hidden _FilterDatabase defined name for this sheet to match.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Espèces!_FilterDatabase") {
        $n.RefersTo = "=Espèces!`$A`$1:`$M`$13"
    }
}

# Make "Espèces" the active sheet with D8 selected (matches the saved view).
$ws.Activate()
$ws.Range("D8").Select()
